$d = $word.ActiveDocument

# Objetivos (PT) body: swap to 'Programa resumido' PT body text
$r = $d.Paragraphs(6).Range
$r.Find.Execute('Proporcionar aos discentes os conhecimentos de biologia celular necessários à compreensão das demais disciplinas do curso e a formação do Engenheiro Ambiental.', $false, $false, $false, $false, $false, $true, 1, $false, 'Análise estrutural das células ao microscópio; moléculas orgânicas; organização interna da célula; organelas celulares transdutoras de energia; material genético e mecanismo de divisão celular.', 2) | Out-Null

# Objetivos (EN) italic body: swap to 'Programa resumido' EN body text
$r = $d.Paragraphs(7).Range
$r.Find.Execute('Provide students with the knowledge of cell biology necessary to understand the other subjects of the course and the training of the Environmental Engineer.', $false, $false, $false, $false, $false, $true, 1, $false, 'Organic molecules; internal organization of the cell; cell energy conversion; genetic material and mechanism of cell division.', 2) | Out-Null

# Docente(s) list bullet: swap to 'Objetivos' PT text
$r = $d.Paragraphs(9).Range
$r.Find.Execute('6712818 - Mauricio Lamano Ferreira', $false, $false, $false, $false, $false, $true, 1, $false, 'Proporcionar aos discentes os conhecimentos de biologia celular necessários à compreensão das demais disciplinas do curso e a formação do Engenheiro Ambiental.', 2) | Out-Null

# Programa resumido (PT) body: swap to 'Programa' PT multiline body text
$r = $d.Paragraphs(11).Range
$r.Find.Execute('Análise estrutural das células ao microscópio; moléculas orgânicas; organização interna da célula; organelas celulares transdutoras de energia; material genético e mecanismo de divisão celular.', $false, $false, $false, $false, $false, $true, 1, $false, '- Estrutura celular e história evolutiva: microrganismos procarióticos e^leucarióticos e suas relações evolutivas dentre os domínios Bacteria, Archaea e^lEukarya.^l- Análise estrutural das células ao microscópio: microscopia ótica e eletrônica.^l- Organização interna da célula: estrutura e função da membrana plasmática; compartimentos intracelulares e seleção de proteínas; tráfico de vesículas (via de exocitose e endocitose).^l- Núcleo e organização do material genético: estrutura e função^l- Ciclo celular e divisão celular: mitose e meiose.^l- Organelas celulares transdutoras de energia: mitocôndria e cloroplasto.', 2) | Out-Null

# Programa resumido (EN) italic body: swap to 'Objetivos' EN text
$r = $d.Paragraphs(12).Range
$r.Find.Execute('Organic molecules; internal organization of the cell; cell energy conversion; genetic material and mechanism of cell division.', $false, $false, $false, $false, $false, $true, 1, $false, 'Provide students with the knowledge of cell biology necessary to understand the other subjects of the course and the training of the Environmental Engineer.', 2) | Out-Null

# Programa (PT) multiline body: swap to 'Metodo' value text
$r = $d.Paragraphs(14).Range
$r.Find.Execute('- Estrutura celular e história evolutiva: microrganismos procarióticos e^leucarióticos e suas relações evolutivas dentre os domínios Bacteria, Archaea e^lEukarya.^l- Análise estrutural das células ao microscópio: microscopia ótica e eletrônica.^l- Organização interna da célula: estrutura e função da membrana plasmática; compartimentos intracelulares e seleção de proteínas; tráfico de vesículas (via de exocitose e endocitose).^l- Núcleo e organização do material genético: estrutura e função^l- Ciclo celular e divisão celular: mitose e meiose.^l- Organelas celulares transdutoras de energia: mitocôndria e cloroplasto.', $false, $false, $false, $false, $false, $true, 1, $false, 'O método de avaliação será composto por avaliação teórica, apresentação escrita e oral.', 2) | Out-Null

# Avaliacao paragraph - 'Norma de recuperacao' value -> bibliography multiline list
$r = $d.Paragraphs(17).Range
$r.Find.Execute('Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 => 5,0 Aprovado', $false, $false, $false, $false, $false, $true, 1, $false, '-Alberts, B. et al. Biologia Molecular da Célula, 5ed. Artmed Editora Ltda, 2010.^l-Cooper, G.M.; Robert, E.H. A célula: uma abordagem molecular. Artmed Editora Ltda, 3ª Edição, 2007.^l-Wasserman, S.A.; Monorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora, 8ª Edição, 2010.^l-Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.^l-Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14 Edição, 2016. ^l-De Roberts, E.M.F.; Hibs, J. Bases da biologia celular e molecular. Editora Guanabara Koogan, 2006.^l-Taiz, L.; Zeiger, E. Plant Physiology. Mass. Sinauer Associates, 2006.', 2) | Out-Null

# Avaliacao paragraph - 'Criterio' value -> old 'Norma de recuperacao' value
$r = $d.Paragraphs(17).Range
$r.Find.Execute('Para o cálculo da nota final (NF) será adotada a média ponderada de provas e atividades.', $false, $false, $false, $false, $false, $true, 1, $false, 'Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 => 5,0 Aprovado', 2) | Out-Null

# Avaliacao paragraph - 'Metodo' value -> old 'Criterio' value
$r = $d.Paragraphs(17).Range
$r.Find.Execute('O método de avaliação será composto por avaliação teórica, apresentação escrita e oral.', $false, $false, $false, $false, $false, $true, 1, $false, 'Para o cálculo da nota final (NF) será adotada a média ponderada de provas e atividades.', 2) | Out-Null

# Bibliography list paragraph -> Docente(s) list bullet text
$r = $d.Paragraphs(19).Range
$r.Find.Execute('-Alberts, B. et al. Biologia Molecular da Célula, 5ed. Artmed Editora Ltda, 2010.^l-Cooper, G.M.; Robert, E.H. A célula: uma abordagem molecular. Artmed Editora Ltda, 3ª Edição, 2007.^l-Wasserman, S.A.; Monorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora, 8ª Edição, 2010.^l-Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.^l-Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14 Edição, 2016. ^l-De Roberts, E.M.F.; Hibs, J. Bases da biologia celular e molecular. Editora Guanabara Koogan, 2006.^l-Taiz, L.; Zeiger, E. Plant Physiology. Mass. Sinauer Associates, 2006.', $false, $false, $false, $false, $false, $true, 1, $false, '6712818 - Mauricio Lamano Ferreira', 2) | Out-Null
